# Corrected excel sheets for application fix issues

$wb = $excel.ActiveWorkbook

function Copy-CellFormat {
    # Copies only the formatting (number format / alignment / etc.) of
    # $srcRange onto $dstRange, leaving $dstRange's value untouched.
    param($srcRange, $dstRange)
    $srcRange.Copy()
    $dstRange.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $excel.CutCopyMode = $false
}

# ---------------------------------------------------------------------------
# Sheet "Summary"
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("B2").Value = 848.21
$wsSummary.Range("E2").NumberFormat = "#,##0.00"
$wsSummary.Range("E2").Value = 4151.79
$wsSummary.Range("F2").Value = 858.48

$wsSummary.Range("A3").Value = 181.28
$wsSummary.Range("E3").Value = 129.49
$wsSummary.Range("F3").Value = 41.52

# ---------------------------------------------------------------------------
# Sheet "Repayment schedule"
# ---------------------------------------------------------------------------
$wsRepay = $wb.Worksheets.Item("Repayment schedule")

# New cell P2, formatted like the rest of row 2 (e.g. O2), left blank.
Copy-CellFormat $wsRepay.Range("O2") $wsRepay.Range("P2")

# Row 3
$wsRepay.Range("F3").Value = 848.21
$wsRepay.Range("G3").NumberFormat = "#,##0.00"
$wsRepay.Range("G3").Value = 4151.79
$wsRepay.Range("K3").Value = 900
$wsRepay.Range("L3").Value = 900
Copy-CellFormat $wsRepay.Range("N3") $wsRepay.Range("O3")
$wsRepay.Range("O3").Value = 0

# Row 4
$wsRepay.Range("F4").Value = 858.48
$wsRepay.Range("G4").NumberFormat = "#,##0.00"
$wsRepay.Range("G4").Value = 3293.31
$wsRepay.Range("H4").Value = 41.52
$wsRepay.Range("K4").Value = 900
Copy-CellFormat $wsRepay.Range("N4") $wsRepay.Range("O4")
$wsRepay.Range("O4").Value = 0
$wsRepay.Range("P4").Value = 900

# Row 5
$wsRepay.Range("F5").Value = 858.72
$wsRepay.Range("G5").NumberFormat = "#,##0.00"
$wsRepay.Range("G5").Value = 2434.59
$wsRepay.Range("H5").Value = 41.28
$wsRepay.Range("K5").Value = 900
Copy-CellFormat $wsRepay.Range("N5") $wsRepay.Range("O5")
$wsRepay.Range("O5").Value = 0
$wsRepay.Range("P5").Value = 900

# Row 6
$wsRepay.Range("F6").Value = 875.65
$wsRepay.Range("G6").NumberFormat = "#,##0.00"
$wsRepay.Range("G6").Value = 1558.94
$wsRepay.Range("H6").Value = 24.35
$wsRepay.Range("K6").Value = 900
Copy-CellFormat $wsRepay.Range("N6") $wsRepay.Range("O6")
$wsRepay.Range("O6").Value = 0
$wsRepay.Range("P6").Value = 900

# Row 7
$wsRepay.Range("F7").Value = 884.41
$wsRepay.Range("G7").Value = 674.53
$wsRepay.Range("H7").Value = 15.59
$wsRepay.Range("K7").Value = 900
Copy-CellFormat $wsRepay.Range("N7") $wsRepay.Range("O7")
$wsRepay.Range("O7").Value = 0
$wsRepay.Range("P7").Value = 900

# Row 8
$wsRepay.Range("F8").Value = 674.53
$wsRepay.Range("H8").Value = 6.75
$wsRepay.Range("K8").Value = 681.28
Copy-CellFormat $wsRepay.Range("N8") $wsRepay.Range("O8")
$wsRepay.Range("O8").Value = 0
$wsRepay.Range("P8").Value = 681.28

# Update the saved selection for this sheet.
$wsRepay.Activate()
$wsRepay.Range("A9:XFD9").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "Transactions"
# ---------------------------------------------------------------------------
$wsTxn = $wb.Worksheets.Item("Transactions")

$wsTxn.Range("A2").Value = 89
$wsTxn.Range("E2").Value = 900
$wsTxn.Range("F2").Value = 848.21
$wsTxn.Range("J2").NumberFormat = "#,##0.00"
$wsTxn.Range("J2").Value = 4151.79

$wsTxn.Range("A3").Value = 88

# Update the saved selection for this sheet (also re-activates/tab-selects it).
$wsTxn.Activate()
$wsTxn.Range("A2:L3").Select() | Out-Null
